$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CON row (row 2) values for columns B:E (meanEMG legmaxROM data)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 14.336693754819585
$ws.Range("C2").Value = 39.339980946029236
$ws.Range("D2").Value = 50.765568178273703
$ws.Range("E2").Value = 42.885043917306568

$ws.Range("B3").Value = 30.333193684649491
$ws.Range("C3").Value = 57.342222431918422
$ws.Range("D3").Value = 67.092777807206431
$ws.Range("E3").Value = 42.86379260784507

# Reflect the updated selection left by the author after editing this range
$ws.Range("B1:E3").Select()
